$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.176.05"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.855.12"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6903"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07784"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3045"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08061"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").Value = "1.857.07"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.189"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.32"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "29.180.20"
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.733"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007805"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "234.85"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").Value = "2.111.25"
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.474"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.14"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.967"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1423"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.95%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.398"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.505"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.480"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05207"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.183"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("E36").Value = "  -2.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.011"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.674"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01850"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.683"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9405"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.96%  "
$ws.Range("D42").Value = "1.096.06"
$ws.Range("E42").Value = "  +4.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.969"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.45"
$ws.Range("D45").ClearFormats()
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.40"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.796"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("D49").Value = "2.007.32"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.162"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.002"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.46%  "
